# Auto-generated Excel COM-interop script to apply numeric updates
# to the Lamia Profits workbook (per-sheet Leve profit recalculation).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1067.3544
$ws.Range("I15").Value = 1067.3544
$ws.Range("K15").Value = 3202.0632
$ws.Range("M15").Value = -3033.0632

$ws.Range("H98").Value = 232627.27
$ws.Range("I98").Value = 969.56525
$ws.Range("K98").Value = 969.56525
$ws.Range("M98").Value = 528.43475

$ws.Range("H106").Value = 3391.1304
$ws.Range("I106").Value = 2419.5
$ws.Range("K106").Value = 2419.5
$ws.Range("M106").Value = -1788.5

$ws.Range("H107").Value = 479.76923
$ws.Range("I107").Value = 520.63635
$ws.Range("K107").Value = 520.63635
$ws.Range("M107").Value = 1399.36365

$ws.Range("H115").Value = 504
$ws.Range("I115").Value = 504
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 1512
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = 55
$ws.Range("N115").ClearContents()

$ws.Range("H118").Value = 812.875
$ws.Range("I118").Value = 592.1667
$ws.Range("K118").Value = 1776.5001
$ws.Range("M118").Value = -119.5001

$ws.Range("H122").Value = 232627.27
$ws.Range("I122").Value = 969.56525
$ws.Range("K122").Value = 2908.69575
$ws.Range("M122").Value = -458.6957499999999

$ws.Range("H132").Value = 1235.8148
$ws.Range("I132").Value = 1271.8077
$ws.Range("K132").Value = 3815.4231
$ws.Range("M132").Value = -1285.4231

$ws.Range("H137").Value = 12197985
$ws.Range("J137").Value = 3167.049
$ws.Range("L137").Value = 9501.147000000001
$ws.Range("N137").Value = -14601.147

$ws.Range("H140").Value = 69977
$ws.Range("J140").Value = 69977
$ws.Range("L140").Value = 69977
$ws.Range("N140").Value = -80337

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7528.8
$ws.Range("I32").Value = 7528.8
$ws.Range("K32").Value = 7528.8
$ws.Range("M32").Value = -7241.8

$ws.Range("H45").Value = 4305
$ws.Range("I45").Value = 3763.6667
$ws.Range("J45").Value = 5117
$ws.Range("K45").Value = 3763.6667
$ws.Range("L45").Value = 5117
$ws.Range("M45").Value = -3386.6667
$ws.Range("N45").Value = -5871

$ws.Range("H61").Value = 6882.591
$ws.Range("I61").Value = 6529.8237
$ws.Range("K61").Value = 6529.8237
$ws.Range("M61").Value = -6317.8237

$ws.Range("H63").Value = 7965.6
$ws.Range("I63").Value = 5414.8335
$ws.Range("J63").Value = 9666.111000000001
$ws.Range("K63").Value = 5414.8335
$ws.Range("L63").Value = 9666.111000000001
$ws.Range("M63").Value = -4728.8335
$ws.Range("N63").Value = -11038.111

$ws.Range("H66").Value = 7965.6
$ws.Range("I66").Value = 5414.8335
$ws.Range("J66").Value = 9666.111000000001
$ws.Range("K66").Value = 27074.1675
$ws.Range("L66").Value = 48330.55500000001
$ws.Range("M66").Value = -23642.1675
$ws.Range("N66").Value = -55194.55500000001

$ws.Range("H136").Value = 6882.591
$ws.Range("I136").Value = 6529.8237
$ws.Range("K136").Value = 19589.4711
$ws.Range("M136").Value = -17039.4711

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 8748
$ws.Range("I20").Value = 8783.286
$ws.Range("J20").Value = 8665.666999999999
$ws.Range("K20").Value = 8783.286
$ws.Range("L20").Value = 8665.666999999999
$ws.Range("M20").Value = -8536.286
$ws.Range("N20").Value = -9159.666999999999

$ws.Range("H86").Value = 4751.1924
$ws.Range("I86").Value = 3741.2
$ws.Range("J86").Value = 8117.8335
$ws.Range("K86").Value = 3741.2
$ws.Range("L86").Value = 8117.8335
$ws.Range("M86").Value = -2618.2
$ws.Range("N86").Value = -10363.8335

$ws.Range("H89").Value = 4751.1924
$ws.Range("I89").Value = 3741.2
$ws.Range("J89").Value = 8117.8335
$ws.Range("K89").Value = 18706
$ws.Range("L89").Value = 40589.1675
$ws.Range("M89").Value = -13090
$ws.Range("N89").Value = -51821.1675

$ws.Range("H94").Value = 1057.7322
$ws.Range("I94").Value = 1146.0834
$ws.Range("J94").Value = 527.625
$ws.Range("K94").Value = 1146.0834
$ws.Range("L94").Value = 527.625
$ws.Range("M94").Value = -695.0834
$ws.Range("N94").Value = -1429.625

$ws.Range("H134").Value = 3802.5557
$ws.Range("I134").Value = 3186.1667
$ws.Range("J134").Value = 5035.3335
$ws.Range("K134").Value = 9558.500100000001
$ws.Range("L134").Value = 15106.0005
$ws.Range("M134").Value = -7023.500100000001
$ws.Range("N134").Value = -20176.0005

$ws.Range("H140").Value = 67039.55499999999
$ws.Range("J140").Value = 67039.55499999999
$ws.Range("L140").Value = 67039.55499999999
$ws.Range("N140").Value = -77399.55499999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 6246.625
$ws.Range("I22").Value = 1493.5
$ws.Range("K22").Value = 1493.5
$ws.Range("M22").Value = -1143.5

$ws.Range("H31").Value = 41271.895
$ws.Range("I31").Value = 2773.4614
$ws.Range("K31").Value = 2773.4614
$ws.Range("M31").Value = -2478.4614

$ws.Range("H34").Value = 41271.895
$ws.Range("I34").Value = 2773.4614
$ws.Range("K34").Value = 2773.4614
$ws.Range("M34").Value = -2571.4614

$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()

$ws.Range("H50").Value = 29900
$ws.Range("J50").Value = 29900
$ws.Range("L50").Value = 29900
$ws.Range("N50").Value = -31150

$ws.Range("H59").Value = 33199.6
$ws.Range("J59").Value = 50000
$ws.Range("L59").Value = 50000
$ws.Range("N59").Value = -52290

$ws.Range("H60").Value = 36126.453
$ws.Range("J60").Value = 36739.1
$ws.Range("L60").Value = 36739.1
$ws.Range("N60").Value = -37761.1

$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws.Range("H132").Value = 3969.1738
$ws.Range("I132").Value = 3556.1428
$ws.Range("K132").Value = 10668.4284
$ws.Range("M132").Value = -8138.428400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 51.416668
$ws.Range("I33").Value = 36.625
$ws.Range("J33").Value = 81
$ws.Range("K33").Value = 219.75
$ws.Range("L33").Value = 486
$ws.Range("M33").Value = 63.25
$ws.Range("N33").Value = -1052

$ws.Range("H122").Value = 4090.7273
$ws.Range("J122").Value = 4404.9
$ws.Range("L122").Value = 39644.1
$ws.Range("N122").Value = -44544.1

$ws.Range("H132").Value = 5893.45
$ws.Range("J132").Value = 6828.3335
$ws.Range("L132").Value = 61455.0015
$ws.Range("N132").Value = -66515.0015

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 23080.555
$ws.Range("I99").Value = 14032.143
$ws.Range("J99").Value = 54750
$ws.Range("K99").Value = 14032.143
$ws.Range("L99").Value = 54750
$ws.Range("M99").Value = -11786.143
$ws.Range("N99").Value = -59242

$ws.Range("H126").Value = 4133.4585
$ws.Range("I126").Value = 2978.5833
$ws.Range("J126").Value = 5288.3335
$ws.Range("K126").Value = 8935.749899999999
$ws.Range("L126").Value = 15865.0005
$ws.Range("M126").Value = -6465.749899999999
$ws.Range("N126").Value = -20805.0005

$ws.Range("H132").Value = 7562.96
$ws.Range("I132").Value = 4929.9443
$ws.Range("K132").Value = 14789.8329
$ws.Range("M132").Value = -12259.8329

$ws.Range("H141").Value = 67097.5
$ws.Range("J141").Value = 67097.5
$ws.Range("L141").Value = 67097.5
$ws.Range("N141").Value = -77457.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6048.081
$ws.Range("I7").Value = 4567.0293
$ws.Range("J7").Value = 22833.334
$ws.Range("K7").Value = 4567.0293
$ws.Range("L7").Value = 22833.334
$ws.Range("M7").Value = -4455.0293
$ws.Range("N7").Value = -23057.334

$ws.Range("H93").Value = 2087.9285
$ws.Range("I93").Value = 2112.476
$ws.Range("J93").Value = 2014.2858
$ws.Range("K93").Value = 2112.476
$ws.Range("L93").Value = 2014.2858
$ws.Range("M93").Value = -864.4760000000001
$ws.Range("N93").Value = -4510.2858

$ws.Range("H122").Value = 104516.48
$ws.Range("I122").Value = 109411.945
$ws.Range("J122").Value = 11502.5
$ws.Range("K122").Value = 328235.835
$ws.Range("L122").Value = 34507.5
$ws.Range("M122").Value = -325785.835
$ws.Range("N122").Value = -39407.5

$ws.Range("H126").Value = 6048.081
$ws.Range("I126").Value = 4567.0293
$ws.Range("J126").Value = 22833.334
$ws.Range("K126").Value = 13701.0879
$ws.Range("L126").Value = 68500.00199999999
$ws.Range("M126").Value = -11231.0879
$ws.Range("N126").Value = -73440.00199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1897.4375
$ws.Range("I122").Value = 1457
$ws.Range("K122").Value = 4371
$ws.Range("M122").Value = -1921

$ws.Range("H126").Value = 3015.4167
$ws.Range("I126").Value = 1855.2778
$ws.Range("J126").Value = 6495.8335
$ws.Range("K126").Value = 5565.8334
$ws.Range("L126").Value = 19487.5005
$ws.Range("M126").Value = -3095.8334
$ws.Range("N126").Value = -24427.5005
